# contratos-11-2020.xlsx — "fix: fixed formatting when scrapping floating
# point numbers"
#
# 1) A handful of "Razon social" / "Nombre Fantasia" values used a comma as
#    an ad-hoc separator between multiple people/parties (e.g.
#    "MONTICO, RICARDO", "ALBIZZATTI, PABLO MARTIN Y FULINI, SERGIO RUBEN").
#    Those stray commas are replaced with periods (and, where a value had a
#    redundant trailing "." or internal "S.H.", those periods are dropped)
#    so the comma is no longer confused with the decimal separator used
#    elsewhere in the sheet.
# 2) The "Importe" column (H) held amounts formatted with the
#    Spanish/Argentine convention — "." as the thousands separator and ","
#    as the decimal separator (e.g. "358.500,00"). They are rewritten to the
#    plain/invariant numeric-text form "358500.00" (no thousands separator,
#    "." as decimal separator) while remaining plain text cells (same as
#    before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Razon social / Nombre Fantasia: stray "," used as a separator -> "." ---
$nameFixes = @{
    "MONTICO, RICARDO" = "MONTICO. RICARDO"
    "ALBIZZATTI, PABLO MARTIN Y FULINI, SERGIO RUBEN" = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "SCHAB DARIO, PEROTTI XAVIER, BENINCA MATIAS S.H." = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
    "URUMAT SOCIEDAD SIMPLE DE BONASEGLA CATALINA, BONASEGLA LUCIANA Y BONASEGLA SILVIO." = "URUMAT SOCIEDAD SIMPLE DE BONASEGLA CATALINA. BONASEGLA LUCIANA Y BONASEGLA SILVIO"
    "DENING BLANCO, CRISTIAN DAVID" = "DENING BLANCO. CRISTIAN DAVID"
    "PARRAVICINI VIRGINIA VANINA, VIRGINIA VANINA" = "PARRAVICINI VIRGINIA VANINA. VIRGINIA VANINA"
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in 5, 6) {
        $cell = $ws.Cells.Item($r, $c)
        $old = $cell.Text
        if ($nameFixes.ContainsKey($old)) {
            $cell.Value = $nameFixes[$old]
        }
    }
}

# --- 2) Importe (column H): "1.234,56" (es-AR) -> "1234.56" (invariant) ---
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Text
    if ($old -ne $null -and $old -ne "") {
        $new = $old -replace '\.', '' -replace ',', '.'
        if ($new -ne $old) {
            # Force the value to stay plain text (it would otherwise be
            # auto-coerced to a number by the numeric-looking string),
            # then drop back to the workbook's default/general style so no
            # visible formatting changes.
            $cell.NumberFormat = "@"
            $cell.Value = $new
            $cell.Style = "Normal"
        }
    }
}
